$wb = $excel.ActiveWorkbook

# --- 1. Update "Metadata" sheet: Last Updated timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value2 = "05 Nov 2025, 10:51 AM"

# --- 2. Update "Stock List" sheet: insert a new top row (CAPTRU-RE1), ---
#         shift all other data rows down by one, dropping the last row.
$ws = $wb.Worksheets.Item("Stock List")

# Shift rows 76..3 down from rows 75..2 (work bottom-up so we don't clobber
# data before it has been copied).
for ($r = 76; $r -ge 3; $r--) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $ws.Cells.Item($r - 1, $c).Value2
    }
}

# Populate the new row 2 with the CAPTRU-RE1 entry.
$ws.Cells.Item(2, 1).Value2 = "📋"
$ws.Cells.Item(2, 2).Value2 = "CAPTRU-RE1"
$ws.Cells.Item(2, 3).Value2 = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value2 = 5.67
$ws.Cells.Item(2, 5).Value2 = -11.9565
$ws.Cells.Item(2, 6).Value2 = "N/A"
$ws.Cells.Item(2, 7).Value2 = "N/A"
$ws.Cells.Item(2, 8).Value2 = 0
